$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert 10 new "costs_type" rows right after the existing block (old row 60 -> new row 70) ---
$ws.Rows("60:69").Insert()

$ws.Range("A60").Value = "costs_type"
$ws.Range("B60").Value = "price_growth"
$ws.Range("C60").Value = "Price Growth"
$ws.Range("D60").Value = "Рост цен"

$ws.Range("A61").Value = "costs_type"
$ws.Range("B61").Value = "price_growth_comp"
$ws.Range("C61").Value = "Price Growth Competitors"
$ws.Range("D61").Value = "Рост цен конкурентов"

$ws.Range("A62").Value = "costs_type"
$ws.Range("B62").Value = "di_growth"
$ws.Range("C62").Value = "Distribution Growth"
$ws.Range("D62").Value = "Рост дистрибьюции"

$ws.Range("A63").Value = "costs_type"
$ws.Range("B63").Value = "di_growth_comp"
$ws.Range("C63").Value = "Distribution Growth Competitors"
$ws.Range("D63").Value = "Рост дистрибьюции конкурентов"

$ws.Range("A64").Value = "costs_type"
$ws.Range("B64").Value = "beauty"
$ws.Range("C64").Value = "Beauty"

$ws.Range("A65").Value = "costs_type"
$ws.Range("B65").Value = "gift_pack"
$ws.Range("D65").Value = "Подарочные упаковки"

$ws.Range("A66").Value = "costs_type"
$ws.Range("B66").Value = "posm"
$ws.Range("C66").Value = "Promotion"
$ws.Range("D66").Value = "Промоушн"

$ws.Range("A67").Value = "costs_type"
$ws.Range("B67").Value = "pr"

$ws.Range("A68").Value = "costs_type"
$ws.Range("B68").Value = "sampling"
$ws.Range("C68").Value = "Sampling"
$ws.Range("D68").Value = "Сэмплинг"

$ws.Range("A69").Value = "costs_type"
$ws.Range("B69").Value = "tailor_made"
$ws.Range("C69").Value = "Tailor Made"
$ws.Range("D69").Value = "На заказ"

$ws.Range("C67").Value = "PR"
$ws.Range("D67").Value = "Пиар"

$ws.Range("D64").Value = "Бьюти-консультанты"

$ws.Range("C65").Value = "Gift Packs"

# --- Apply the (re-)applied font to column B of all "costs_type" rows (25-69), which Excel
#     records as a brand-new (duplicate) font entry + cellXf ---
$ws.Range("B25:B69").Font.ThemeColor = 1

# --- View cosmetics: zoom to 100% and move the selection ---
$excel.ActiveWindow.Zoom = 100
$ws.Range("C66").Select()
